# Rename the three embedded logo pictures (the first-page footer's
# Pearson logo, the default footer's Pearson logo, and the first-page
# header's BTec logo) to match the target revision:
#
#   footer1.xml (first-page footer, docPr id="3") : image1.png -> image2.png
#   footer2.xml (default footer,    docPr id="2") : image1.png -> image2.png
#   header1.xml (first-page header, docPr id="1") : image2.jpg -> image1.jpg
#
# Inline pictures don't expose a settable "Name" directly, so each
# picture is converted to a floating shape (which does expose .Name),
# renamed, then converted back to an inline picture in place.

$d   = $word.ActiveDocument
$sec = $d.Sections(1)

# --- First-page footer: Pearson logo (wp:docPr id="3") ---
$ftrFirst = $sec.Footers(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -ge 1) {
    $shp = $ftrFirst.Range.InlineShapes(1)
    $floating = $shp.ConvertToShape()
    $floating.Name = "image2.png"
    $floating.ConvertToInlineShape()
}

# --- Default footer: Pearson logo (wp:docPr id="2") ---
$ftrDefault = $sec.Footers(1)
if ($ftrDefault.Exists -and $ftrDefault.Range.InlineShapes.Count -ge 1) {
    $shp = $ftrDefault.Range.InlineShapes(1)
    $floating = $shp.ConvertToShape()
    $floating.Name = "image2.png"
    $floating.ConvertToInlineShape()
}

# --- First-page header: BTec logo (wp:docPr id="1") ---
$hdrFirst = $sec.Headers(2)
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -ge 1) {
    $shp = $hdrFirst.Range.InlineShapes(1)
    $floating = $shp.ConvertToShape()
    $floating.Name = "image1.jpg"
    $floating.ConvertToInlineShape()
}
